$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Responsable " column header in I1
$ws.Range("I1").Value = "Responsable "

# Fill I2:I227 with the responsible person's name for every record
$ws.Range("I2:I227").Value = "Geovani Hernández Gómez"

# Match the vertical-center alignment style used by the rest of the data columns
$ws.Range("I2:I227").VerticalAlignment = -4108

# Move the active selection to I2, matching the saved selection state
[void]$ws.Range("I2").Select()
